# Apply updated cryptocurrency price/volume figures to the "cryptos" sheet.
# D-column (Price) values are forced to remain plain text (matching the
# original inline-string cell type) even when they look like plain numbers
# (e.g. "245.70"), by briefly switching the cell to a Text number format
# while the value is entered, then restoring General.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.375.34"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.884.89"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.70%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.70"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.689"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.68"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +5.42%  "
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.37"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.158.80"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.763"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +4.20%  "
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.876.78"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.513.16"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.08"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.68"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.82"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.96"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("E24").Value = "  +9.11%  "
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.16"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -5.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.72"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.54"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.32"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("E30").Value = "  -1.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128.45"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.72"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +9.56%  "
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("E34").Value = "  -3.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.91"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -6.68%  "
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.846"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0699"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +7.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.40"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0219"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.31"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -5.95%  "
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.308.07"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("E46").Value = "  -3.24%  "
$ws.Range("E47").Value = "  +7.11%  "
$ws.Range("E48").Value = "  -2.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.73"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.16"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.25"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -4.89%  "
